$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Txn Date Min" column (G) values for rows 2-5
$ws.Range("G2").Value = 43104
$ws.Range("G3").Value = 43101
$ws.Range("G4").Value = 43101
$ws.Range("G5").Value = 43102
